$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C45").Value = 'Opaline vases with bronze base'
$ws.Range("C47").Value = 'Opaline flower pineapple pot'
$ws.Range("C48").Value = 'Green Ottoman sarm'
$ws.Range("C51").Value = 'Moroccan necklace, old frame, black thread center'
$ws.Range("C52").Value = 'Moroccan necklace, old frame, blue thread center'
$ws.Range("C53").Value = 'Iranian vase, porcelain, blue, inlaid pottery'
$ws.Range("C55").Value = 'Transparent decorated Bohemian vase'
$ws.Range("C58").Value = 'Navy graffiti pair, large size'
$ws.Range("C59").Value = 'Blue French lantern pair'
$ws.Range("C60").Value = 'Blue Bohemian lantern pair'
$ws.Range("C61").Value = 'Turquoise Bohemian lantern large'
$ws.Range("C62").Value = 'Blue and white Czech step boss pair'
$ws.Range("C74").Value = 'Moroccan Luban frame, various types of stones, large'
$ws.Range("C75").Value = 'Moroccan Luban frame, various types of stones, small'
$ws.Range("C76").Value = 'Sea lion fang candlestick with silver'
$ws.Range("C77").Value = 'French Bohemian red box with silver elephant handle'
$ws.Range("C79").Value = 'Green and black opaline walnuts'
$ws.Range("C83").Value = 'White transparent opaline cup'
$ws.Range("C86").Value = 'French Limoges navy box'
$ws.Range("C95").Value = 'Manuscript plaque of supplications'
$ws.Range("C96").Value = 'Ottoman calligraphy plaque (The head of wisdom is the fear of God) Ottoman calligrapher Rasan'
$ws.Range("C98").Value = 'Hameed al-Ghasi calligraphy plaque'
$ws.Range("C103").Value = 'Iranian saffron sprinklers'
$ws.Range("C104").Value = 'Bohemian box with French bronze'
$ws.Range("C108").Value = 'Crystal box with bronze'
$ws.Range("C125").Value = 'Large frame Moroccan necklace with different stones'
$ws.Range("C128").Value = 'Syrian seashell caskets'
$ws.Range("C130").Value = 'Copper applique nut'
$ws.Range("C132").Value = 'Ottoman sirma silver frame'
$ws.Range("C97").Value = 'Moroccan wood and brass table'
$ws.Range("C46").Value = 'Opaline vases with a picture of a gazelle head'
$ws.Range("C49").Value = 'Red Bohemian graffiti pair / 2'
$ws.Range("C50").Value = 'Red Bohemian sprinklers pair / 2'
$ws.Range("C54").Value = 'Blue and green French maison plates pair, pictured'
$ws.Range("C56").Value = 'Syrian consul, 5 drawers, wooden'
$ws.Range("C57").Value = 'Ottoman sarm, large, burgundy'
$ws.Range("C63").Value = 'Transparent base cup, Italian blue'
$ws.Range("C64").Value = 'Consul with mirror, Syrian seashell'
$ws.Range("C65").Value = 'Opaline walnut, illustrated, rose, beige, vertical, brown'
$ws.Range("C66").Value = 'Opaline walnut, illustrated, French'
$ws.Range("C67").Value = 'Opaline walnut, illustrated, medium-sized, white and yellow rose'
$ws.Range("C68").Value = 'Opaline single, green color, vase'
$ws.Range("C69").Value = 'Opaline walnut, white, gilded and green flower'
$ws.Range("C70").Value = 'Opaline walnut, rose yellow'
$ws.Range("C71").Value = 'Opaline walnut, large-sized, illustrated, maroon walnut'
$ws.Range("C72").Value = 'Opaline vase, large-sized, illustrated and signed'
$ws.Range("C73").Value = '3-door cabinet, Syrian-Arab workshop, Ottoman era, seashell, 120 years old'
$ws.Range("C78").Value = 'Lemon-colored walnuts'
$ws.Range("C80").Value = '3-piece set, parvatin walnuts with plate'
$ws.Range("C81").Value = 'Bohemian green gilded box'
$ws.Range("C82").Value = 'Navy walnut Oriental vases'
$ws.Range("C84").Value = 'French navy blue box illustrated with bronze'
$ws.Range("C85").Value = 'French Limoges egg walnut navy blue with bronze'
$ws.Range("C87").Value = '3-piece set of fat and meat vases'
$ws.Range("C88").Value = 'Opaline fat and meat hookah'
$ws.Range("C89").Value = 'Small frame colored stones crystal'
$ws.Range("C90").Value = 'Old Syrian mother of pearl two-door cabinet'
$ws.Range("C91").Value = 'Indian bone chairs pair Maharaja'
$ws.Range("C92").Value = 'Desk with chair with drawer unit Old Syrian mother of pearl'
$ws.Range("C93").Value = 'Old Iranian wool carpet Hazrat Suleiman'
$ws.Range("C94").Value = 'Ayat al-Kursi plaque with a reed'
$ws.Range("C99").Value = 'Italian marble columns pair with bronze'
$ws.Range("C100").Value = 'Large old Syrian mother of pearl mirror'
$ws.Range("C101").Value = 'French bronze camel sculptural walnut'
$ws.Range("C102").Value = 'Saffron vases pair with white steppe ute'
$ws.Range("C105").Value = 'Blue lanterns'
$ws.Range("C106").Value = 'Green vase with a picture of a gazelle head'
$ws.Range("C107").Value = '3-piece set of French opaline for the Islamic market antique'
$ws.Range("C109").Value = 'Illustrated bronze French box'
$ws.Range("C110").Value = 'Black bohemian box with bronze'
$ws.Range("C111").Value = 'Illustrated black round box'
$ws.Range("C112").Value = 'Old Syrian display cabinet with shell'
$ws.Range("C113").Value = 'Wool carpet walnuts with Iranian silk inlay / 2'
$ws.Range("C114").Value = '2 rooster figurines'
$ws.Range("C115").Value = '3-piece bronze ballerina set figurines'
$ws.Range("C116").Value = 'Porcelain rooster figurines / 2'
$ws.Range("C117").Value = 'White elephant figurine with wooden base'
$ws.Range("C118").Value = 'English porcelain tiger figurine'
$ws.Range("C119").Value = 'Glass turtle figurine'
$ws.Range("C120").Value = 'Gemstone elephant figurine (good stone)'
$ws.Range("C121").Value = 'Vertical turquoise lizard figurine'
$ws.Range("C122").Value = 'Horizontal turquoise lizard figurine'
$ws.Range("C123").Value = 'Turquoise crocodile figurine'
$ws.Range("C124").Value = 'Medium Syrian shell mirrors'
$ws.Range("C126").Value = 'Large Ottoman sirma with silver and silk threads'
$ws.Range("C127").Value = 'Black sirmat figurines with gold trim'
$ws.Range("C129").Value = 'Back shell boxes Turtle'
$ws.Range("C131").Value = 'French yellow and pink vase nut'
$ws.Range("C133").Value = 'Moroccan vertical large painting'
$ws.Range("C134").Value = 'Moroccan vertical large painting'
$ws.Range("C135").Value = 'Moroccan vertical large painting'
$ws.Range("C136").Value = 'Moroccan vertical large painting'
$ws.Range("C137").Value = 'Moroccan vertical large painting'
$ws.Range("C138").Value = 'Moroccan vertical large painting'
$ws.Range("C139").Value = 'Moroccan painting in the middle of a picture of an old man'
$ws.Range("C140").Value = 'Moroccan painting in the middle of a picture of a desert'
$ws.Range("C141").Value = 'Moroccan painting in the middle of a picture of a desert'
$ws.Range("C142").Value = 'Moroccan painting in the middle of a picture of a woman'
$ws.Range("C143").Value = 'Moroccan painting in the middle of a square picture of a desert ... woman gold frame'
$ws.Range("C144").Value = 'Moroccan painting in the middle of a picture of a woman gold frame'
$ws.Range("C145").Value = 'Moroccan painting in the picture of a horse vertical'
$ws.Range("C146").Value = 'Moroccan painting in the picture of a horse vertical'
$ws.Range("C147").Value = 'Moroccan painting in the picture of a horse vertical'
$ws.Range("C148").Value = 'Moroccan painting in the picture of a horse square'
$ws.Range("C149").Value = 'Moroccan painting in the picture of a horse vertical'
$ws.Range("C150").Value = 'Bedouin Arabic Painting Gold Frame'
$ws.Range("C151").Value = 'Moroccan Souq Painting'
$ws.Range("C152").Value = 'Moroccan Woman with Pottery Painting'
$ws.Range("C153").Value = 'Bird Visit Painting'
$ws.Range("C154").Value = 'Large Moroccan Painting in Lily'
$ws.Range("C155").Value = 'Kadoddle Painting'
$ws.Range("C156").Value = 'Moroccan Square Woman with Jar Gold and Red Frame'
$ws.Range("C157").Value = 'Western Souq Painting in Scalet'
$ws.Range("C158").Value = 'Moroccan Wood and Brass Table'
$ws.Range("C159").Value = 'Black Frame Body Painting with Silver'
$ws.Range("C160").Value = 'Colored Moroccan Painting'
$ws.Range("C161").Value = 'Moroccan Woman Gold Frame Painting'
$ws.Range("C162").Value = 'Biro Shell with Small Mirror'
$ws.Range("C163").Value = 'European Cabinet 2 Doors Wood with French Picture'
$ws.Range("C164").Value = 'European Cabinet 3 Doors Wood English'
$ws.Range("C165").Value = '4 Piece Set Opaline White'
$ws.Range("C166").Value = 'White Opaline Perfume Bottles Stainless Steel 2 Pieces'
$ws.Range("C167").Value = 'White Opaline Lantern Single Small'
$ws.Range("C168").Value = 'White Opaline Lantern Single Medium'
$ws.Range("C169").Value = 'White Decorated Perfume Box'
$ws.Range("C170").Value = 'White Decorated Perfume Vase'
$ws.Range("C171").Value = 'Small Rose Perfume Vase'
$ws.Range("C172").Value = 'White Perfume Box Decorated with Gold'
$ws.Range("C173").Value = 'White Opaline Cup'
$ws.Range("C174").Value = 'Blue Opaline Chandeliers'
$ws.Range("C175").Value = 'White and Blue'
$ws.Range("C176").Value = 'Tharya Step Pot White and Blue'
$ws.Range("C177").Value = 'French Brass Candlestick Nut with Navy Blue Picture'
$ws.Range("C178").Value = 'French Pure Brass Candlestick'

$win = $excel.ActiveWindow
try { $win.ScrollRow = 158 } catch {}
try { $win.ScrollColumn = 1 } catch {}
[void]$ws.Range("C182").Select()

Write-Host "done"
